$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "test"
$ws.Range("B3").Value = "test"
$ws.Range("C3").Value = "test.png"
$ws.Range("D3").Value = "test"
$ws.Range("E3").Value = 22
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
